$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their text formatting so Excel
# does not auto-convert numeric-looking strings (e.g. "1.004") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.800.55"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.633.23"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "214.98"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "0.5076"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.2576"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "0.06415"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "20.26"
$ws.Range("E10").Value = "  +4.02%  "
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.647.95"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "4.247"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").Value = "1.859.47"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "0.5575"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "0.0₅7639"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "63.17"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "25.808.43"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "4.360"
$ws.Range("D21").Value = "191.71"
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "9.903"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -6.48%  "
$ws.Range("D26").Value = "138.97"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "0.1225"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "1.240"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "0.04940"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "3.288"
$ws.Range("E32").Value = "  +1.77%  "
$ws.Range("D33").Value = "3.245"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("D34").Value = "1.568"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").Value = "2.386"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "0.8995"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "0.5548"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").Value = "1.130.44"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "5.448"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "98.88"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "0.7989"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "0.0₈111"
$ws.Range("E45").Value = "  -4.70%  "
$ws.Range("D46").Value = "55.52"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "0.4259"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").Value = "7.770"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").Value = "0.05025"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").Value = "0.9962"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  +0.27%  "
